# Update computed profit-margin figures (columns H:N) on the Balmung_Profits sheets
# per the scheduled runner refresh. Generated from the row-level cell diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 31
$ws.Range("H31").Value = 541.6
$ws.Range("I31").Value = 176
$ws.Range("J31").Value = 2004
$ws.Range("K31").Value = 528
$ws.Range("L31").Value = 6012
$ws.Range("M31").Value = -298
$ws.Range("N31").Value = -6472

# Row 40
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("N40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("M40").ClearContents()

# Row 64
$ws.Range("H64").Value = 10676.6
$ws.Range("I64").Value = 6483.3335
$ws.Range("J64").Value = 16966.5
$ws.Range("K64").Value = 6483.3335
$ws.Range("L64").Value = 16966.5
$ws.Range("M64").Value = -6235.3335
$ws.Range("N64").Value = -17462.5

# Row 67
$ws.Range("H67").Value = 10676.6
$ws.Range("I67").Value = 6483.3335
$ws.Range("J67").Value = 16966.5
$ws.Range("K67").Value = 6483.3335
$ws.Range("L67").Value = 16966.5
$ws.Range("M67").Value = -5625.3335
$ws.Range("N67").Value = -18682.5

# Row 80
$ws.Range("H80").Value = 17857608
$ws.Range("I80").Value = 311.5
$ws.Range("J80").Value = 25000526
$ws.Range("K80").Value = 934.5
$ws.Range("L80").Value = 75001578
$ws.Range("M80").Value = 63.5
$ws.Range("N80").Value = -75003574

# Row 83
$ws.Range("H83").Value = 17857608
$ws.Range("I83").Value = 311.5
$ws.Range("J83").Value = 25000526
$ws.Range("K83").Value = 2803.5
$ws.Range("L83").Value = 225004734
$ws.Range("M83").Value = 2188.5
$ws.Range("N83").Value = -225014718

# Row 137
$ws.Range("H137").Value = 1288247.8
$ws.Range("I137").Value = 5947.028
$ws.Range("J137").Value = 2387362.5
$ws.Range("K137").Value = 17841.084
$ws.Range("L137").Value = 7162087.5
$ws.Range("M137").Value = -15291.084
$ws.Range("N137").Value = -7167187.5

# Row 138
$ws.Range("H138").Value = 26394.2
$ws.Range("I138").Value = 59000
$ws.Range("J138").Value = 4657
$ws.Range("K138").Value = 177000
$ws.Range("L138").Value = 13971
$ws.Range("M138").Value = -171860
$ws.Range("N138").Value = -24251

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2167.2134
$ws.Range("I32").Value = 1134.8422
$ws.Range("J32").Value = 8202.615
$ws.Range("K32").Value = 1134.8422
$ws.Range("L32").Value = 8202.615
$ws.Range("M32").Value = -847.8422
$ws.Range("N32").Value = -8776.615

# Row 122
$ws.Range("H122").Value = 1837.9062
$ws.Range("I122").Value = 1215.7407
$ws.Range("J122").Value = 5197.6
$ws.Range("K122").Value = 3647.2221
$ws.Range("L122").Value = 15592.8
$ws.Range("M122").Value = -1197.2221
$ws.Range("N122").Value = -20492.8

# Row 125
$ws.Range("H125").Value = 124994.5
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 124994.5
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 124994.5
$ws.Range("N125").Value = -134834.5

$ws = $wb.Worksheets.Item("BSM")
# Row 130
$ws.Range("H130").Value = 54631.668
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 54631.668
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 54631.668
$ws.Range("N130").Value = -64671.668

$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 5015000
$ws.Range("I6").Value = 5015000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 5015000
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -5014887

# Row 7
$ws.Range("H7").Value = 76.44444
$ws.Range("I7").Value = 94
$ws.Range("J7").Value = 41.333332
$ws.Range("K7").Value = 94
$ws.Range("L7").Value = 41.333332
$ws.Range("M7").Value = 19
$ws.Range("N7").Value = -267.333332

# Row 25
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("N25").Value = 0
$ws.Range("L25").ClearContents()

# Row 31
$ws.Range("H31").Value = 4538.4688
$ws.Range("I31").Value = 3244.9375
$ws.Range("J31").Value = 4797.175
$ws.Range("K31").Value = 3244.9375
$ws.Range("L31").Value = 4797.175
$ws.Range("M31").Value = -2949.9375
$ws.Range("N31").Value = -5387.175

# Row 34
$ws.Range("H34").Value = 4538.4688
$ws.Range("I34").Value = 3244.9375
$ws.Range("J34").Value = 4797.175
$ws.Range("K34").Value = 3244.9375
$ws.Range("L34").Value = 4797.175
$ws.Range("M34").Value = -3042.9375
$ws.Range("N34").Value = -5201.175

# Row 41
$ws.Range("H41").Value = 5029.5
$ws.Range("I41").Value = 5029.5
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 5029.5
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -4601.5

# Row 50
$ws.Range("H50").Value = 15000
$ws.Range("I50").Value = 15000
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 15000
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -14375

# Row 58
$ws.Range("H58").Value = 2578.5
$ws.Range("I58").Value = 2255.5557
$ws.Range("J58").Value = 3159.8
$ws.Range("K58").Value = 2255.5557
$ws.Range("L58").Value = 3159.8
$ws.Range("M58").Value = -2052.5557
$ws.Range("N58").Value = -3565.8

# Row 62
$ws.Range("H62").Value = 8450
$ws.Range("I62").Value = 8450
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 8450
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = -7826
$ws.Range("M62").ClearContents()

# Row 65
$ws.Range("H65").Value = 8450
$ws.Range("I65").Value = 8450
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 42250
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = -39130
$ws.Range("M65").ClearContents()

# Row 136
$ws.Range("H136").Value = 2578.5
$ws.Range("I136").Value = 2255.5557
$ws.Range("J136").Value = 3159.8
$ws.Range("K136").Value = 6766.6671
$ws.Range("L136").Value = 9479.400000000001
$ws.Range("M136").Value = -4216.6671
$ws.Range("N136").Value = -14579.4

$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 1634.4286
$ws.Range("I7").Value = 298.33334
$ws.Range("J7").Value = 2636.5
$ws.Range("K7").Value = 895.0000200000001
$ws.Range("L7").Value = 7909.5
$ws.Range("M7").Value = -783.0000200000001
$ws.Range("N7").Value = -8133.5

# Row 34
$ws.Range("H34").Value = 2552.7144
$ws.Range("I34").Value = 543.8
$ws.Range("J34").Value = 7575
$ws.Range("K34").Value = 1631.4
$ws.Range("L34").Value = 22725
$ws.Range("M34").Value = -1547.4
$ws.Range("N34").Value = -22893

# Row 39
$ws.Range("H39").Value = 4914.923
$ws.Range("I39").Value = 1999.5
$ws.Range("J39").Value = 5445
$ws.Range("K39").Value = 5998.5
$ws.Range("L39").Value = 16335
$ws.Range("M39").Value = -5704.5
$ws.Range("N39").Value = -16923

# Row 114
$ws.Range("H114").Value = 1811.5
$ws.Range("I114").Value = 1742.25
$ws.Range("J114").Value = 1950
$ws.Range("K114").Value = 5226.75
$ws.Range("L114").Value = 5850
$ws.Range("M114").Value = -1972.75
$ws.Range("N114").Value = -12358

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 20925822
$ws.Range("I80").Value = 38907.562
$ws.Range("J80").Value = 62699652
$ws.Range("K80").Value = 38907.562
$ws.Range("L80").Value = 62699652
$ws.Range("M80").Value = -37909.562
$ws.Range("N80").Value = -62701648

# Row 83
$ws.Range("H83").Value = 20925822
$ws.Range("I83").Value = 38907.562
$ws.Range("J83").Value = 62699652
$ws.Range("K83").Value = 194537.81
$ws.Range("L83").Value = 313498260
$ws.Range("M83").Value = -189545.81
$ws.Range("N83").Value = -313508244

# Row 132
$ws.Range("H132").Value = 2851028.5
$ws.Range("I132").Value = 2372.1667
$ws.Range("J132").Value = 8548341
$ws.Range("K132").Value = 7116.500100000001
$ws.Range("L132").Value = 25645023
$ws.Range("M132").Value = -4586.500100000001
$ws.Range("N132").Value = -25650083

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 6465.091
$ws.Range("I16").Value = 5840.2
$ws.Range("J16").Value = 6985.8335
$ws.Range("K16").Value = 5840.2
$ws.Range("L16").Value = 6985.8335
$ws.Range("M16").Value = -5670.2
$ws.Range("N16").Value = -7325.8335

# Row 68
$ws.Range("H68").Value = 5034.3125
$ws.Range("I68").Value = 5265.1
$ws.Range("J68").Value = 4649.6665
$ws.Range("K68").Value = 5265.1
$ws.Range("L68").Value = 4649.6665
$ws.Range("M68").Value = -4516.1
$ws.Range("N68").Value = -6147.6665

# Row 69
$ws.Range("H69").Value = 60000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 60000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 60000
$ws.Range("N69").Value = -61622

# Row 71
$ws.Range("H71").Value = 5034.3125
$ws.Range("I71").Value = 5265.1
$ws.Range("J71").Value = 4649.6665
$ws.Range("K71").Value = 26325.5
$ws.Range("L71").Value = 23248.3325
$ws.Range("M71").Value = -22581.5
$ws.Range("N71").Value = -30736.3325

# Row 72
$ws.Range("H72").Value = 60000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 60000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 180000
$ws.Range("N72").Value = -188112

# Row 130
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("N130").Value = 0
$ws.Range("L130").ClearContents()
